# Update the weekly units report:
#  - refresh the "Report Generated On" timestamp
#  - insert a new billable line item (PLA-TAG) at Point 01 (row 17),
#    pushing the remaining line items down by one row
#  - zero out all pricing ("Pricing" column) values and the grand TOTAL
#  - bump the "Total Line Items" count
#  - keep the "Total Billed Amount" summary in sync (now 0)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- header / summary cells -------------------------------------------------
$ws.Range("D5").Value = "Report Generated On: 08/18/2025 09:49 PM"
$ws.Range("C8").Value = 0
$ws.Range("C9").Value = 36

# --- insert a new detail row at row 17 (shifts 17..51 down to 18..52) ------
$ws.Rows.Item(17).EntireRow.Insert()

# Copy the formatting (fill / font / number format) of the row that is now
# two rows below (a row that already carries the "odd" banding style) onto
# the freshly inserted, still-blank row 17 so the banding pattern stays
# consistent with the rest of the table.
$ws.Range("A20:H20").Copy()
$ws.Range("A17:H17").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- populate the new line item ---------------------------------------------
$ws.Range("A17").Value = "Point 01"
$ws.Range("B17").Value = "PLA-TAG"
$ws.Range("C17").Value = "Inst"
$ws.Range("D17").Value = "PLA,Tag Pole"
$ws.Range("E17").Value = "EA"
$ws.Range("F17").Value = 5
$ws.Range("H17").Value = 0

# --- zero out every pricing value in the (now 16..52) detail rows ----------
for ($r = 16; $r -le 52; $r++) {
    $ws.Cells.Item($r, 8).Value = 0
}

# --- nudge the used range out to column I, matching the original layout ----
$ws.Range("I52").NumberFormat = "General"
